$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout:  N=VERMOEGEN  O=EL-BEZUG  P=SH-BEZUG
# Target layout:    N=EGID  O=EWID  P=VERMOEGEN  Q=STEUERBARESEINKOMMEN  R=EL-BEZUG  S=SH-BEZUG  T=AMOUNT

# Step 1: insert two blank columns before N -> VERMOEGEN/EL-BEZUG/SH-BEZUG shift to P/Q/R
$ws.Columns("N:O").Insert()

# Step 2: insert one blank column before the (shifted) EL-BEZUG column (now Q) to make room
#         for STEUERBARESEINKOMMEN -> EL-BEZUG/SH-BEZUG shift to R/S
$ws.Columns("Q").Insert()

# New headers
$ws.Range("N1").Value2 = "EGID"
$ws.Range("O1").Value2 = "EWID"
$ws.Range("Q1").Value2 = "STEUERBARESEINKOMMEN"
$ws.Range("S1").Copy($ws.Range("T1"))
$ws.Range("T1").Value2 = "AMOUNT"

# Row 2 data
$ws.Range("N2").Value2 = 11111
$ws.Range("Q2").Value2 = 12345
$ws.Range("T2").Value2 = 100

# Row 3 data
$ws.Range("N3").Value2 = 22
$ws.Range("Q3").Value2 = 99999
$ws.Range("T3").Value2 = 2000

# Row 4 data
$ws.Range("N4").Value2 = 9
$ws.Range("O4").Value2 = 122
$ws.Range("Q4").Value2 = 400000
$ws.Range("T4").Value2 = -50

# Update selection to match target
$ws.Range("T5").Select()
